$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("llama3.1-8b-instruct")

# Two new data rows are being added. Insert two blank rows above the
# existing "h2o (mc4000)" row (currently row 5) so it shifts down to row 7,
# opening up rows 5 and 6 for the new entries.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Insert()

# Row 3: was "max_fused (mc4000)" -> becomes "h2o (mc1000)"
$ws.Range("A3").Value = "preds_ns10_ws200_mc1000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse"
$ws.Range("B3").Value = 21156.49535999999

# Row 4: was "sum_fused" -> becomes "max_fused" (value unchanged)
$ws.Range("A4").Value = "preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse"
$ws.Range("B4").Value = 53752.95488000003

# Row 5 (new): ws32 max_fused
$ws.Range("A5").Value = "preds_ns10_ws32_mc4000_ea1.0_snks0_hopf_True_type_max_fused_rerun_lenNone_gblFalse"
$ws.Range("B5").Value = 26220.95360000002
$ws.Range("A5").Borders.LineStyle = 1

# Row 6 (new): h2o (mc4000) moved up from the old row 5 position
$ws.Range("A6").Value = "preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_h2o_rerun_lenNone_gblFalse"
$ws.Range("B6").Value = 84562.57535999996
$ws.Range("A6").Borders.LineStyle = 1

# Row 7: already holds the shifted-down "h2o (mc4000)" row's old content -
# overwrite it with the relocated "sum_fused" entry.
$ws.Range("A7").Value = "preds_ns10_ws200_mc4000_ea1.0_snks0_hopf_True_type_sum_fused_rerun_lenNone_gblFalse"
$ws.Range("B7").Value = 53752.95488000003
